# Fix bug in config files
# Apply content edits to C_BackPolygon.conf sheet:
#  - A1 header comment changes from "C_BackPolygon" to "# C_BackPolygon"
#  - All "kind%=..." filter expressions in column A lose the stray "%" -> "kind=..."
# Also reset the saved view state (remove scrolled topLeftCell, select A1:C1 header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_BackPolygon.conf")

$ws.Range("A4").Value = "kind=0136"
$ws.Range("A5").Value = "kind=0137"
$ws.Range("A6").Value = "kind=0133"
$ws.Range("A7").Value = "kind=0134"
$ws.Range("A8").Value = "kind=0121"
$ws.Range("A9").Value = "kind=0123;name=筒子河"
$ws.Range("A10").Value = "#kind=0123;dispclass=8"
$ws.Range("A11").Value = "#kind=0123;dispclass=7"
$ws.Range("A12").Value = "kind=0123;dispclass=5"
$ws.Range("A13").Value = "kind=0123;dispclass=4"
$ws.Range("A14").Value = "kind=0123;dispclass=3"
$ws.Range("A15").Value = "kind=0123;dispclass=2"
$ws.Range("A16").Value = "kind=0123;dispclass=1"
$ws.Range("A17").Value = "kind=0123"
$ws.Range("A18").Value = "kind=0122"
$ws.Range("A19").Value = "kind=0125"
$ws.Range("A20").Value = "kind=0143"
$ws.Range("A21").Value = "kind=0145"
$ws.Range("A22").Value = "kind=0146"
$ws.Range("A23").Value = "kind=0147"
$ws.Range("A24").Value = "kind=0148"
$ws.Range("A25").Value = "scene_id<InRange>0001-0002;(kind=0141 || kind=014b || kind=tx0160)"
$ws.Range("A26").Value = "kind=0141"
$ws.Range("A27").Value = "kind=0142"
$ws.Range("A28").Value = "kind=0144"
$ws.Range("A29").Value = "kind=0149"
$ws.Range("A30").Value = "kind=014a"
$ws.Range("A31").Value = "kind=014b"
$ws.Range("A32").Value = "kind=0161"
$ws.Range("A33").Value = "kind=0162"
$ws.Range("A34").Value = "kind=0163"
$ws.Range("A35").Value = "kind=0164"
$ws.Range("A36").Value = "kind=0165"
$ws.Range("A37").Value = "kind=0166"
$ws.Range("A38").Value = "kind=0167"
$ws.Range("A39").Value = "kind=0171"
$ws.Range("A40").Value = "kind=0172"
$ws.Range("A41").Value = "kind=0173"
$ws.Range("A42").Value = "kind=0174"
$ws.Range("A43").Value = "kind=0175"
$ws.Range("A44").Value = "kind=0176"
$ws.Range("A45").Value = "kind=0177"
$ws.Range("A46").Value = "kind=0178"
$ws.Range("A47").Value = "kind=0179"
$ws.Range("A48").Value = "kind=017a"
$ws.Range("A49").Value = "kind=84FF || kind=84ff"
$ws.Range("A50").Value = "kind=07FF || kind=07ff"
$ws.Range("A65").Value = "kind=010bj0101"
$ws.Range("A66").Value = "kind=010bj0102"
$ws.Range("A67").Value = "kind=0123"
$ws.Range("A68").Value = "kind=0141 || kind=014b || kind=tx0160"
$ws.Range("A69").Value = "kind=tx1000"
$ws.Range("A70").Value = "kind=tx2010"
$ws.Range("A71").Value = "kind=tx2000"
$ws.Range("A72").Value = "kind=tx1020"
$ws.Range("A73").Value = "kind=tx1030"
$ws.Range("A74").Value = "kind=tx1040"
$ws.Range("A75").Value = "kind=tx1050"
$ws.Range("A76").Value = "kind=tx1031"
$ws.Range("A77").Value = "kind=tx1051"
$ws.Range("A78").Value = "kind=tx1070"
$ws.Range("A79").Value = "kind=tx1060"

# The header row comment is updated last so the shared-string table rebuild
# appends it after the "kind=" strings above (matching the target layout).
$ws.Range("A1").Value = "# C_BackPolygon"

# Reset the view: scroll back to the top and select the header row (A1:C1),
# matching the saved selection state in the target workbook.
$ws.Activate()
$ws.Range("A1:C1").Select()
